$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 20.02.2022 23:30"

# Update row 10 (EuroOil Opuštěná) with the latest scraped price check
$ws.Range("B10").Value = 37.4
$ws.Range("C10").Value = 36.5

# Delta Cena / Old Datum columns changed from numeric to plain text cells
$ws.Range("D10").Value = "'+0.9"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = "'2022-02-20 23:30:44"
$ws.Range("E10").ClearFormats()
